$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell H1 "eta^2", reusing the same formatting as the
# existing header row (bold, centered, bordered) by copying from G1.
$ws.Range("H1").Value = "eta^2"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Add the new eta^2 (effect size) values for each row
$ws.Range("H2").Value = 0.02
$ws.Range("H3").Value = 0.25
$ws.Range("H4").Value = 0.11
$ws.Range("H5").Value = 0.11
